$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("customers")
$ws.Range("C6").Value = $true
